# Fix column A references: strip the stray trailing "16" that was
# appended to every "<Book> <Chapter>:<Verse>" reference string, e.g.
# "2 Timothy 1:116" -> "2 Timothy 1:1". Any cell that does not end in
# "16" (e.g. row 13 "2 Timothy 1:18") is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value2

    if ($value -ne $null -and $value.ToString().EndsWith("16")) {
        $fixed = $value.ToString().Substring(0, $value.ToString().Length - 2)
        $cell.Value = $fixed
    }
}
